$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New loading_percent values for rows 2-25 (A2:A25 = 0..23), columns B:N.
# Column order: B,C,D,E,F,G,H,I,J,K,L,M,N
$data = @(
    @(18.89028610029903, 6.345714281256337, 4.378109198514716, 10.4293559780028, 51.6151704710622, 0, 7.344005520526261, 37.34093666912207, 9.982248629945566, 16.93002746922939, 11.80697136171512, 18.258407811372, 25.60714598775363),
    @(18.79215134195066, 6.219701808648395, 4.37883125161381, 10.44275322461274, 51.64967152507997, 0, 7.344005520526261, 37.39697471182208, 9.994377578097199, 16.86433294254147, 11.8228106336086, 18.26322245136709, 25.66099979129025),
    @(18.73593723290198, 6.14299805608556, 4.379658938110845, 10.45162027778671, 51.67953640869833, 0, 7.344005520526261, 37.43669907214399, 10.00223420293586, 16.82754616970421, 11.83390621630253, 18.26899284731115, 25.6959881038659),
    @(18.71406519095668, 6.111956571874271, 4.380093312408247, 10.45539529493567, 51.69388899746383, 0, 7.344005520526261, 37.45422327083315, 10.00553910023714, 16.81345958079614, 11.83877274431747, 18.27205337791766, 25.71073014396532),
    @(18.71049644733641, 6.106816602165165, 4.380171316279412, 10.45603190733881, 51.69640401798777, 0, 7.344005520526261, 37.45721382894674, 10.00609412166712, 16.81117545139066, 11.8396016763351, 18.27260443792619, 25.71320730430454),
    @(18.73563804088588, 6.142578482080134, 4.379664402555584, 10.45167053403384, 51.67972113749506, 0, 7.344005520526261, 37.43693000055282, 10.00227835543138, 16.82735251690923, 11.83397045059724, 18.26903125011816, 25.69618495942849),
    @(18.85562267614436, 6.302155016003608, 4.378278657951051, 10.43384253889242, 51.62526453038225, 0, 7.344005520526261, 37.35915468090465, 9.986345907708712, 16.906646528625, 11.81214862646497, 18.2594847135892, 25.62531634961661),
    @(19.12205066372562, 6.61831654004658, 4.378591483282995, 10.4039503287904, 51.58735985544737, 0, 7.344005520526261, 37.24885406759974, 9.958336656971497, 17.0897481350354, 11.78020842412788, 18.26302487900165, 25.50155787286486),
    @(19.33544772260257, 6.849874373355082, 4.380642648630129, 10.38505352178071, 51.60148158303893, 0, 7.344005520526261, 37.1935892917588, 9.939710114260906, 17.24028403889941, 11.76333084699563, 18.27909436816203, 25.41985938350777),
    @(19.43605230865415, 6.954518858086259, 4.381965514925262, 10.37711716149404, 51.61700064356369, 0, 7.344005520526261, 37.1740489406305, 9.931656033431729, 17.31204920798101, 11.75707723805251, 18.2893028902167, 25.38468607534728),
    @(19.47462864632362, 6.994003238533258, 4.382521985479924, 10.37420634434596, 51.62418216077611, 0, 7.344005520526261, 37.1674547736263, 9.928666127576248, 17.33967956459617, 11.75491335339758, 18.29358269226974, 25.37165256437991),
    @(19.46629966647536, 6.985506550202476, 4.382399678309266, 10.37482904358711, 51.622577518626, 0, 7.344005520526261, 37.16883912852362, 9.929307393299334, 17.33370895644383, 11.75537031116328, 18.29264259120612, 25.3744468591895),
    @(19.43921655890106, 6.957770348093876, 4.382010186216693, 10.37687579486663, 51.61756534168421, 0, 7.344005520526261, 37.17349029455585, 9.931408851214679, 17.31431337225077, 11.75689512439107, 18.28964671033751, 25.38360807465348),
    @(19.42268898381569, 6.940761376448271, 4.381778827027668, 10.37814178586211, 51.61466505642226, 0, 7.344005520526261, 37.17644414619363, 9.932703860200631, 17.30249165099556, 11.75785569495334, 18.2878654838093, 25.38925679461286),
    @(19.32894162471225, 6.843017979629059, 4.380563995418235, 10.38558542765721, 51.60065014585351, 0, 7.344005520526261, 37.19497899520998, 9.940244878921968, 17.23565860850166, 11.76376814415648, 18.27848531293033, 25.4221980732376),
    @(19.27231497225583, 6.782847173293567, 4.379918201730164, 10.39032060297958, 51.59437995950197, 0, 7.344005520526261, 37.2077839098758, 9.944978220494379, 17.19548765523165, 11.76775959471199, 18.27347137202423, 25.44291621532),
    @(19.2400787342735, 6.74817576272838, 4.379583478270386, 10.39310628544875, 51.59162981060214, 0, 7.344005520526261, 37.21567605181363, 9.947740189851654, 17.17269277518724, 11.77018944882207, 18.27086047903552, 25.45502024692764),
    @(19.22922228505535, 6.736427178039764, 4.379476468981538, 10.39406015303718, 51.59084580845474, 0, 7.344005520526261, 37.21843871854597, 9.948682134051987, 17.16502866665985, 11.771035200105, 18.2700234416843, 25.45915068526267),
    @(19.2783086167574, 6.789259269502487, 4.379983150927555, 10.38981010753275, 51.59495882784313, 0, 7.344005520526261, 37.20636625133235, 9.944470264604496, 17.19973192570812, 11.76732082543852, 18.27397687984118, 25.44069133176807),
    @(19.44715873553915, 6.965921331083602, 4.382123086632999, 10.3762720524517, 51.6190021561811, 0, 7.344005520526261, 37.17210227796018, 9.930789975964711, 17.31999813940711, 11.7564417119112, 18.29051545880882, 25.38090944982922),
    @(19.56029380934061, 7.080533147663258, 4.383845078157124, 10.36797485325247, 51.6423191997156, 0, 7.344005520526261, 37.15440282571188, 9.922198710642984, 17.40123903334488, 11.75052165636802, 18.30373637189117, 25.34350448604346),
    @(19.49966639350919, 7.019453654845566, 4.382896601078717, 10.37235296083552, 51.6291799554829, 0, 7.344005520526261, 37.16341987563317, 9.926752136152531, 17.35764382910229, 11.75357260354872, 18.29646039071248, 25.36331596268077),
    @(19.27559789372868, 6.786360601774534, 4.379953673521015, 10.39004070512416, 51.59469445890898, 0, 7.344005520526261, 37.20700552243997, 9.944699784636466, 17.1978121553698, 11.76751877211005, 18.27374749325335, 25.44169660066271),
    @(19.04677360148703, 6.532724494468948, 4.378185079583277, 10.41149692908399, 51.5902391650743, 0, 7.344005520526261, 37.27416923844829, 9.965569712335915, 17.0373394231509, 11.78768994376049, 18.25969360619345, 25.5334140532648)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    for ($j = 0; $j -lt $data[$i].Length; $j++) {
        $colNum = $j + 2
        $ws.Cells.Item($rowNum, $colNum).Value = $data[$i][$j]
    }
}
